$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 217, shifting existing rows 217:222 down to 218:223
$ws.Rows.Item(217).Insert()

# Fill the new row 217 with data (copy fixed columns from its neighbor, then set unique values)
$ws.Cells.Item(217, 1).Value = 4
$ws.Cells.Item(217, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(217, 3).Value = "Los Lagos"
$ws.Cells.Item(217, 4).Value = 44448
$ws.Cells.Item(217, 5).Value = 10
$ws.Cells.Item(217, 6).Value = 100114001
$ws.Cells.Item(217, 7).Value = "Papa"
$ws.Cells.Item(217, 8).Value = "Asterix"
$ws.Cells.Item(217, 9).Value = "1a (guarda)"
$ws.Cells.Item(217, 10).Value = 300
$ws.Cells.Item(217, 11).Value = 7000
$ws.Cells.Item(217, 12).Value = 7500
$ws.Cells.Item(217, 13).Value = 7250
$ws.Cells.Item(217, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(217, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(217, 16).Value = 290
$ws.Cells.Item(217, 17).Value = 25
$ws.Cells.Item(217, 18).Value = "Hortaliza"
